# Addition of decentralized generators and changes in import
$wb = $excel.ActiveWorkbook

$ppObj1 = $wb.Worksheets.Item("PP Object Type 1")

# Fix the truncated label and add the new decentralized capacity row
$ppObj1.Range("A5").Value = "NL Installed Capacity-RES (+heat)"
$ppObj1.Range("A6").Value = "NL Installed Capacity Decentralized (+heat)"
$ppObj1.Range("B6").Value = "UNITNL"

# Resize column A to fit the new, longer labels
$ppObj1.Columns.Item(1).AutoFit() | Out-Null

# Make "PP Object Type 1" the active sheet / tab, with B7 selected
$ppObj1.Activate()
$ppObj1.Range("B7").Select() | Out-Null
